# Automatic update of files.
# Update the "Förändrad" (Changed) date column (C) for rows 2-23 from
# serial date 45171 (2023-09-02) to 45172 (2023-09-03), preserving the
# existing date number format/style already applied to these cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 23; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45171) {
        $cell.Value2 = 45172
    }
}
